{"js": "// Locate the paragraph that ends the \"Socks in the Dark\" evaluation section,\n// then insert the new \"Choose a solution...\" block right after it (mirroring\n// the structure already used for Problem 1).\nconst searchResults = context.document.body.search(\n  \"Each of these solutions will meet the goals\",\n  { matchCase: true }\n);\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Anchor paragraph not found\");\n}\n\nconst anchorParagraphs = searchResults.items[0].paragraphs;\nanchorParagraphs.load(\"items\");\nawait context.sync();\n\nconst anchor = anchorParagraphs.items[0];\n\n// 1) blank spacer paragraph\nconst blankPara = anchor.insertParagraph(\"\", Word.InsertLocation.after);\nawait context.sync();\n\n// 2) \"Choose a solution...\" heading line\nconst choosePara = blankPara.insertParagraph(\n  \"Choose a solution and develop a plan to implement it:\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\n// 3) the full explanation paragraph\nconst detailPara = choosePara.insertParagraph(\n  \"The first solution goes like this: to guarantee a pair of one color every time, 4 socks are necessary because two socks picked could be of different colors.  Also, 3 socks picked could still be one of each color and the absence of a pair because there are three colors of socks to be possibly chosen.  A fourth sock will guarantee solution because the fourth sock must be one of the three colors already picked.  The second solution uses the same ideology to come up with another guaranteed solution.  This time, the worst case scenario would be that you pick 17 socks up and they are 10 black, 6 brown, and only one white sock.  The next pick would have to be another white sock and would solve the problem.  Thus, 18 socks must be picked to guarantee a pair of each color is picked every time.\",\n  Word.InsertLocation.after\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the paragraph that ends the \"Socks in the Dark\" evaluation section\n# (\"Each of these solutions will meet the goals...\") so we can insert the new\n# \"Choose a solution...\" block right after it, mirroring the layout already\n# used for Problem 1.\n$targetIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n  $t = $d.Paragraphs.Item($i).Range.Text\n  if ($t -like \"Each of these solutions will meet the goals*\") {\n    $targetIndex = $i\n    break\n  }\n}\n\nif ($targetIndex -eq -1) {\n  throw \"Anchor paragraph not found\"\n}\n\n$anchor = $d.Paragraphs.Item($targetIndex)\n\n# 1) blank spacer paragraph (stays empty, like the spacers used elsewhere in\n#    this document)\n$anchor.Range.InsertParagraphAfter()\n$spacer = $d.Paragraphs.Item($targetIndex + 1)\n\n# 2) \"Choose a solution...\" heading line -- insert a fresh paragraph after the\n#    spacer and give it its own text, so the spacer itself stays blank.\n$spacer.Range.InsertParagraphAfter()\n$choosePara = $d.Paragraphs.Item($targetIndex + 2)\n$choosePara.Range.Text = \"Choose a solution and develop a plan to implement it:\"\n\n# 3) the full explanation paragraph -- same pattern: insert a fresh paragraph\n#    after the heading line and set its text.\n$choosePara.Range.InsertParagraphAfter()\n$detailPara = $d.Paragraphs.Item($targetIndex + 3)\n$detailPara.Range.Text = \"The first solution goes like this: to guarantee a pair of one color every time, 4 socks are necessary because two socks picked could be of different colors.  Also, 3 socks picked could still be one of each color and the absence of a pair because there are three colors of socks to be possibly chosen.  A fourth sock will guarantee solution because the fourth sock must be one of the three colors already picked.  The second solution uses the same ideology to come up with another guaranteed solution.  This time, the worst case scenario would be that you pick 17 socks up and they are 10 black, 6 brown, and only one white sock.  The next pick would have to be another white sock and would solve the problem.  Thus, 18 socks must be picked to guarantee a pair of each color is picked every time.\"\n"}
